$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @(2, 3, 5),
    @(2, 6, 217),
    @(2, 8, "kitchens"),
    @(2, 12, "stimuli/img_p3hpc.png"),
    @(2, 13, 72.83333333333333),
    @(2, 14, 52.22222222222222),
    @(2, 15, 62.52777777777777),
    @(2, 16, 36),
    @(2, 17, 6),
    @(2, 18, 6),
    @(2, 19, 6),
    @(3, 3, 5),
    @(3, 6, 218),
    @(3, 8, "kitchens"),
    @(3, 9, "target"),
    @(3, 11, "j"),
    @(3, 12, "stimuli/img_t90e2.png"),
    @(3, 13, 83.0625),
    @(3, 14, 61.96875),
    @(3, 15, 72.515625),
    @(3, 16, 32),
    @(3, 17, 9),
    @(3, 18, 9),
    @(3, 19, 9),
    @(4, 3, 5),
    @(4, 6, 219),
    @(4, 8, "kitchens"),
    @(4, 12, "stimuli/img_7wul8.png"),
    @(4, 13, 43.03030303030303),
    @(4, 14, 25.54545454545455),
    @(4, 15, 34.28787878787879),
    @(4, 16, 33),
    @(4, 17, 1),
    @(4, 18, 1),
    @(4, 19, 1),
    @(5, 3, 5),
    @(5, 6, 220),
    @(5, 8, "kitchens"),
    @(5, 12, "stimuli/img_iyxnj.png"),
    @(5, 13, 75.30555555555556),
    @(5, 14, 54.33333333333334),
    @(5, 15, 64.81944444444444),
    @(5, 16, 36),
    @(5, 17, 6),
    @(5, 18, 6),
    @(5, 19, 6),
    @(6, 3, 5),
    @(6, 6, 221),
    @(6, 8, "kitchens"),
    @(6, 12, "stimuli/img_uwv6y.png"),
    @(6, 13, 78.88888888888889),
    @(6, 14, 59.30555555555556),
    @(6, 15, 69.09722222222223),
    @(6, 16, 36),
    @(6, 17, 8),
    @(6, 18, 8),
    @(6, 19, 8),
    @(7, 3, 5),
    @(7, 6, 222),
    @(7, 8, "kitchens"),
    @(7, 12, "stimuli/img_9mky8.png"),
    @(7, 13, 84.32352941176471),
    @(7, 14, 65.17647058823529),
    @(7, 15, 74.75),
    @(7, 16, 34),
    @(7, 17, 9),
    @(7, 18, 9),
    @(7, 19, 9),
    @(8, 3, 5),
    @(8, 6, 223),
    @(8, 8, "kitchens"),
    @(8, 12, "stimuli/img_ce9vx.png"),
    @(8, 13, 75.9090909090909),
    @(8, 14, 57.12121212121212),
    @(8, 15, 66.51515151515152),
    @(8, 16, 33),
    @(8, 17, 7),
    @(8, 18, 7),
    @(8, 19, 7),
    @(9, 3, 5),
    @(9, 6, 224),
    @(9, 8, "kitchens"),
    @(9, 12, "stimuli/img_ye5sl.png"),
    @(9, 13, 53.2258064516129),
    @(9, 14, 34.45161290322581),
    @(9, 15, 43.83870967741936),
    @(9, 16, 31),
    @(9, 17, 2),
    @(9, 18, 2),
    @(9, 19, 2),
    @(10, 3, 5),
    @(10, 6, 225),
    @(10, 8, "kitchens"),
    @(10, 9, "target"),
    @(10, 11, "j"),
    @(10, 12, "stimuli/img_yeh72.png"),
    @(10, 13, 68.66666666666667),
    @(10, 14, 45.21212121212121),
    @(10, 15, 56.93939393939394),
    @(10, 16, 33),
    @(10, 17, 4),
    @(10, 18, 4),
    @(10, 19, 4),
    @(11, 3, 5),
    @(11, 6, 226),
    @(11, 8, "living_rooms"),
    @(11, 12, "stimuli/img_jpjeg.png"),
    @(11, 13, 90.90697674418605),
    @(11, 14, 74.3953488372093),
    @(11, 15, 82.65116279069767),
    @(12, 3, 5),
    @(12, 6, 227),
    @(12, 8, "kitchens"),
    @(12, 12, "stimuli/img_cnyac.png"),
    @(12, 13, 69.1470588235294),
    @(12, 14, 47.8235294117647),
    @(12, 15, 58.48529411764706),
    @(12, 16, 34),
    @(12, 17, 5),
    @(12, 18, 5),
    @(12, 19, 5),
    @(13, 3, 5),
    @(13, 6, 228),
    @(13, 12, "stimuli/img_kn0we.png"),
    @(13, 13, 80.1590909090909),
    @(13, 14, 56.68181818181818),
    @(13, 15, 68.42045454545455),
    @(13, 16, 44),
    @(13, 17, 7),
    @(13, 18, 7),
    @(13, 19, 7),
    @(14, 3, 5),
    @(14, 6, 229),
    @(14, 8, "kitchens"),
    @(14, 9, "target"),
    @(14, 11, "j"),
    @(14, 12, "stimuli/img_463mq.png"),
    @(14, 13, 51.35294117647059),
    @(14, 14, 30.20588235294118),
    @(14, 15, 40.77941176470588),
    @(14, 16, 34),
    @(14, 17, 2),
    @(14, 18, 2),
    @(14, 19, 2),
    @(15, 3, 5),
    @(15, 6, 230),
    @(15, 8, "kitchens"),
    @(15, 12, "stimuli/img_eatdk.png"),
    @(15, 13, 81.40625),
    @(15, 14, 61.375),
    @(15, 15, 71.390625),
    @(15, 16, 32),
    @(15, 17, 8),
    @(15, 18, 8),
    @(15, 19, 8),
    @(16, 3, 5),
    @(16, 6, 231),
    @(16, 8, "kitchens"),
    @(16, 12, "stimuli/img_inqod.png"),
    @(16, 13, 70.84848484848484),
    @(16, 14, 50.63636363636363),
    @(16, 15, 60.74242424242424),
    @(16, 16, 33),
    @(16, 17, 5),
    @(16, 18, 5),
    @(16, 19, 5),
    @(17, 3, 5),
    @(17, 6, 232),
    @(17, 8, "kitchens"),
    @(17, 9, "target"),
    @(17, 11, "j"),
    @(17, 12, "stimuli/img_aplao.png"),
    @(17, 13, 64.0909090909091),
    @(17, 14, 40.75757575757576),
    @(17, 15, 52.42424242424242),
    @(17, 16, 33),
    @(17, 17, 3),
    @(17, 18, 3),
    @(17, 19, 3),
    @(18, 3, 5),
    @(18, 6, 233),
    @(18, 8, "kitchens"),
    @(18, 12, "stimuli/img_3gm8h.png"),
    @(18, 13, 65.07894736842105),
    @(18, 14, 43.92105263157895),
    @(18, 15, 54.5),
    @(18, 16, 38),
    @(18, 17, 4),
    @(18, 18, 4),
    @(18, 19, 4),
    @(19, 3, 5),
    @(19, 6, 234),
    @(19, 8, "kitchens"),
    @(19, 12, "stimuli/img_a8wvq.png"),
    @(19, 13, 86.25925925925925),
    @(19, 14, 66.25925925925925),
    @(19, 15, 76.25925925925925),
    @(19, 16, 27),
    @(19, 17, 10),
    @(19, 18, 10),
    @(19, 19, 10),
    @(20, 3, 5),
    @(20, 6, 235),
    @(20, 9, "distractor"),
    @(20, 11, "f"),
    @(20, 12, "stimuli/img_53nbn.png"),
    @(20, 13, 73.28888888888889),
    @(20, 14, 51.15555555555556),
    @(20, 15, 62.22222222222223),
    @(20, 16, 45),
    @(20, 17, 6),
    @(20, 18, 6),
    @(20, 19, 6),
    @(21, 3, 5),
    @(21, 6, 236),
    @(21, 8, "kitchens"),
    @(21, 12, "stimuli/img_d8xbu.png"),
    @(21, 13, 91.36363636363636),
    @(21, 14, 73.18181818181819),
    @(21, 15, 82.27272727272728),
    @(21, 16, 33),
    @(21, 17, 10),
    @(21, 18, 10),
    @(21, 19, 10),
    @(22, 3, 5),
    @(22, 6, 237),
    @(22, 8, "kitchens"),
    @(22, 12, "stimuli/img_60242.png"),
    @(22, 13, 78.33333333333333),
    @(22, 14, 57.57575757575758),
    @(22, 15, 67.95454545454545),
    @(22, 16, 33),
    @(22, 17, 7),
    @(22, 18, 7),
    @(22, 19, 7),
    @(23, 3, 5),
    @(23, 6, 238),
    @(23, 8, "kitchens"),
    @(23, 12, "stimuli/img_nyv2b.png"),
    @(23, 13, 11.91176470588235),
    @(23, 14, 6.852941176470588),
    @(23, 15, 9.382352941176471),
    @(23, 16, 34),
    @(23, 17, 1),
    @(23, 18, 1),
    @(23, 19, 1),
    @(24, 3, 5),
    @(24, 6, 239),
    @(24, 8, "bedrooms"),
    @(24, 9, "distractor"),
    @(24, 11, "f"),
    @(24, 12, "stimuli/img_scrdm.png"),
    @(24, 13, 78.675),
    @(24, 14, 57.9),
    @(24, 15, 68.2875),
    @(24, 16, 40),
    @(24, 17, 7),
    @(24, 18, 7),
    @(24, 19, 7),
    @(25, 3, 5),
    @(25, 6, 240),
    @(25, 8, "kitchens"),
    @(25, 12, "stimuli/img_wyl6z.png"),
    @(25, 13, 59.8235294117647),
    @(25, 14, 36.23529411764706),
    @(25, 15, 48.02941176470588),
    @(25, 16, 34),
    @(25, 17, 3),
    @(25, 18, 3),
    @(25, 19, 3),
    @(26, 3, 5),
    @(26, 6, 241),
    @(26, 9, "distractor"),
    @(26, 11, "f"),
    @(26, 12, "stimuli/img_g13d5.png"),
    @(26, 13, 73),
    @(26, 14, 51.51111111111111),
    @(26, 15, 62.25555555555556),
    @(26, 16, 45),
    @(26, 17, 6),
    @(26, 18, 6),
    @(26, 19, 6),
    @(27, 3, 5),
    @(27, 6, 242),
    @(27, 9, "distractor"),
    @(27, 11, "f"),
    @(27, 12, "stimuli/img_x9w7o.png"),
    @(27, 13, 92.38888888888889),
    @(27, 14, 72.94444444444444),
    @(27, 15, 82.66666666666666),
    @(27, 16, 36),
    @(27, 17, 10),
    @(27, 18, 10),
    @(27, 19, 10)
)

foreach ($item in $changes) {
    $r = $item[0]
    $c = $item[1]
    $v = $item[2]
    $ws.Cells.Item($r, $c).Value = $v
}
